# Mise à jour de l'application
# Add the new player "Nathanael Beta" (already in A32) to row 32:
# Âge = 19, Date de naissance = 07/08/2006, Poste = AD.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Âge
$ws.Range("B32").Value = 19

# Date de naissance - reuse the existing date style/format from another
# "Date de naissance" cell (C30) so the new cell matches the column's
# formatting, then write the actual date value (07/08/2006 -> serial 38936).
$ws.Range("C30").Copy()
$ws.Range("C32").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C32").Value = 38936

# Poste
$ws.Range("D32").Value = "AD"

# Scroll/selection bookkeeping to match the author's last view state.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E34").Select()
